# Update countries & provincias Spain
# Applies the 30-Apr-2020 19:52 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 19:52"

# --- Estados Unidos (row 4) -------------------------------------------
$ws.Range("B4").Value = 1076129
$ws.Range("C4").Value = 11935
$ws.Range("E4").Value = 864063
$ws.Range("G4").Value = 725
$ws.Range("H4").Value = 62380

# --- Kazajistan (row 61) ----------------------------------------------
$ws.Range("B61").Value = 3402
$ws.Range("C61").Value = 264
$ws.Range("D61").Value = 866
$ws.Range("E61").Value = 2511

# --- Uzbekistan (row 71) -----------------------------------------------
$ws.Range("D71").Value = 1133
$ws.Range("E71").Value = 875

# --- Principado de Andorra (row 98) ------------------------------------
$ws.Range("B98").Value = 745
$ws.Range("C98").Value = 2
$ws.Range("D98").Value = 468
$ws.Range("E98").Value = 235

# --- Togo / Trinidad yTobago / Cabo Verde reorder + refresh (147-149) -
# Cabo Verde's updated numbers move it ahead of Togo and Trinidad y
# Tobago in the sorted list; Togo and Trinidad y Tobago keep their
# existing figures and simply shift one row down each.
$ws.Range("A147").Value = "Cabo Verde"
$ws.Range("B147").Value = 116
$ws.Range("C147").Value = 3
$ws.Range("D147").Value = 2
$ws.Range("E147").Value = 113
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 1

$ws.Range("A148").Value = "Togo"
$ws.Range("B148").Value = 116
$ws.Range("C148").Value = 7
$ws.Range("D148").Value = 65
$ws.Range("E148").Value = 42
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 2
$ws.Range("H148").Value = 9

$ws.Range("A149").Value = "Trinidad yTobago"
$ws.Range("B149").Value = 116
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 72
$ws.Range("E149").Value = 36
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 8

# --- Aruba (row 153) ----------------------------------------------------
$ws.Range("D153").Value = 79
$ws.Range("E153").Value = 19

# --- Eritrea (row 174) ---------------------------------------------------
$ws.Range("D174").Value = 26
$ws.Range("E174").Value = 13
